# Append 6 new daily rows (2020-04-17 .. 2020-04-22, serials 43938-43943)
# to the bottom of every worksheet's data table (rows 36-41), extending
# the used range from A1:H35 to A1:H41 on each sheet.

$wb = $excel.ActiveWorkbook

# New row data per worksheet (in sheet order), columns: Date, B, C, D, E, F, G, H
$sheetRows = @{
    1 = @(
        @(43938, 16557, 451, 1037, 725, 82, 15, 14576),
        @(43939, 16612, 476, 1140, 748, 85, 15, 15464),
        @(43940, 17897, 490, 1175, 774, 93, 18, 16060),
        @(43941, 18563, 507, 1223, 822, 102, 18, 16509),
        @(43942, 19101, 533, 1280, 867, 109, 18, 17337),
        @(43943, 19413, 556, 1419, 893, 112, 18, 18100)
    )
    2 = @(
        @(43938, 35036.72433608084, 777.696666640341, 1312.252160552951, 1687.794202899496, 1407.331993736343, 387.7155181326793, 10792.21475571486),
        @(43939, 35153.11135296098, 820.8062379618677, 1442.591574764092, 1741.338018991479, 1458.819749604746, 387.7155181326793, 11449.69875016291),
        @(43940, 37872.33529279693, 844.9475979019227, 1486.881666971762, 1801.865811095461, 1596.120431920486, 465.2586217592152, 11890.98305274291),
        @(43941, 39281.67626083642, 874.2621064005608, 1547.622364856566, 1913.609427287428, 1750.583699525695, 465.2586217592152, 12223.42709948522),
        @(43942, 40420.15289868214, 919.0960605749485, 1619.751943594771, 2018.369067467397, 1870.721796551968, 465.2586217592152, 12836.48649971381),
        @(43943, 41080.38470352947, 958.7568661907532, 1795.646881219515, 2078.896859571379, 1922.209552420371, 465.2586217592152, 13401.41925620465)
    )
    3 = @(
        @(43938, 1725, 37, 81, 57, 2, 0, 847),
        @(43939, 55, 25, 103, 23, 3, 0, 888),
        @(43940, 1285, 14, 35, 26, 8, 3, 596),
        @(43941, 666, 17, 48, 48, 9, 0, 449),
        @(43942, 538, 26, 57, 45, 7, 0, 828),
        @(43943, 312, 23, 139, 26, 3, 0, 763)
    )
    4 = @(
        @(43938, 3650.320074877058, 63.80216555585947, 102.4999276806066, 132.6955442279603, 34.32517057893519, 0, 627.1271883980849),
        @(43939, 116.3870168801381, 43.10957132152667, 130.3394142111417, 53.54381609198399, 51.48775586840279, 0, 657.4839944480512),
        @(43940, 2719.223939835953, 24.14135994005493, 44.29009220766951, 60.52779210398191, 137.3006823157408, 77.54310362653587, 441.2843025799983),
        @(43941, 1409.34096803949, 29.31450849863813, 60.74069788480389, 111.7436161919666, 154.4632676052084, 0, 332.4440467423142),
        @(43942, 1138.476637845714, 44.83395417438773, 72.12957873820463, 104.7596401799687, 120.1380970262732, 0, 613.0594002285883),
        @(43943, 660.2318048473288, 39.66080561580453, 175.8949376247446, 60.52779210398191, 51.48775586840279, 0, 564.9327564908368)
    )
    5 = @(
        @(43938, 3052.302420871113, 57.93926385613183, 100.4752377511131, 122.9179778111633, 102.9755117368056, 0, 586.9969716198368),
        @(43939, 2773.820031390637, 58.97389356784847, 107.8147387455269, 115.9340017991654, 96.11047762101855, 0, 612.3189903249306),
        @(43940, 2985.856014870452, 51.04173244468757, 103.2591864041666, 116.3996001999653, 89.24544350523152, 15.50862072530717, 585.3680698317899),
        @(43941, 2939.724433634324, 45.52370731553215, 91.87030555076592, 105.2252385807686, 106.4080287946991, 15.50862072530717, 539.1664918435484),
        @(43942, 1806.749727495671, 41.04031189809338, 81.99994214448529, 92.65408175917236, 99.54299467891208, 15.50862072530717, 534.2797864794076),
        @(43943, 1208.732073489725, 36.2120399100824, 96.6789441333129, 78.22053133437666, 102.9755117368056, 15.50862072530717, 521.840900097958)
    )
}

foreach ($sheetIndex in 1..5) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $newRows = $sheetRows[$sheetIndex]

    $startRow = 36
    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $targetRow = $startRow + $i
        $srcRange = $ws.Range("A35:H35")
        $dstRange = $ws.Range("A${targetRow}:H${targetRow}")

        # Copy row 35 (values + formatting, incl. the date style on col A)
        # down into the new row so the new row inherits identical styling.
        $srcRange.Copy($dstRange)

        # Overwrite the copied values with the actual data for this row.
        $values = $newRows[$i]
        $rowArray = New-Object 'object[,]' 1,8
        for ($col = 0; $col -lt 8; $col++) {
            $rowArray[0, $col] = $values[$col]
        }
        $dstRange.Value2 = $rowArray
    }
}

Write-Output "Appended rows 36-41 to sheets 1-5"
